# issue #5: add legislator_id, name, date into dataframe
#
# The stock ("股票") sheet gains three new trailing columns - date,
# legislator_name, legislator_id - populated with this filing's values
# for every existing data row (the same per-filing metadata that's
# already being stamped onto the other sheets in this scrape).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "紀國棟"
$legislatorId = 918
$filingDate = "2012-04-16"

$headerCol = 8   # H
$nameCol = 9     # I
$idCol = 10      # J

$firstDataRow = 2
$lastDataRow = $ws.UsedRange.Rows.Count

# --- Header row (row 1): reuse the formatting already used by the other
# header cells (bold / centered / bordered) so the new columns look the
# same as name/owner/quantity/... .
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(1, $headerCol).Value = "date"
$ws.Cells.Item(1, $nameCol).Value = "legislator_name"
$ws.Cells.Item(1, $idCol).Value = "legislator_id"

# --- Data rows: reuse the plain data-cell formatting from the existing
# columns, then fill every row with this filing's date / legislator
# name / legislator id.
$dataRange = "H" + $firstDataRow + ":J" + $lastDataRow
$ws.Range("G2").Copy() | Out-Null
$ws.Range($dataRange).PasteSpecial(-4122) | Out-Null

# Force the date column to text so Excel doesn't silently reinterpret
# "2012-04-16" as a date serial number - the source value is a plain
# string, not a real Excel date.
$dateRange = "H" + $firstDataRow + ":H" + $lastDataRow
$ws.Range($dateRange).NumberFormat = "@"

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, $headerCol).Value = $filingDate
    $ws.Cells.Item($r, $nameCol).Value = $legislatorName
    $ws.Cells.Item($r, $idCol).Value = $legislatorId
}
